$d = $word.ActiveDocument

# Locate the "18/05/2020" date text in the document.
$rng = $d.Content
$found = $rng.Find.Execute("18/05/2020", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $dateStart = $rng.Start

    # Range covering just the day-of-month portion ("18") that must change to "20".
    $dayRange = $d.Range($dateStart, $dateStart + 2)
    $dayRange.Text = "20"

    # Toggle bold off/on (net no-op) so the edited "20" becomes its own run,
    # distinct from the untouched "/05/2020" remainder, matching how Word
    # splits a run when only part of its text is directly edited.
    $dayRange.Font.Bold = $false
    $dayRange.Font.Bold = $true
}
